# Applies the TRIMAZKON config update:
# 1. ip_address_list: add new row 14 for "533valeo"
# 2. Settings: update default interface index (B1) and disk-window startup flag (B4)
# 3. projects_bin2 (hidden): rename/update row 3 entry to "533valeo" and drop its note (D3)
# 4. Settings_recources: add new row 33 for system tray startup setting

$wb = $excel.ActiveWorkbook

# --- 1. ip_address_list ---
$wsIp = $wb.Worksheets.Item("ip_address_list")
$wsIp.Range("A14").Value = "533valeo"
$wsIp.Range("B14").Value = "192.168.227.27"
$wsIp.Range("C14").Value = "255.255.255.0"
$wsIp.Range("E14").Value = $false

# --- 2. Settings ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B1").Value = 6
$wsSettings.Range("B4").Value = 0

# --- 3. projects_bin2 (hidden sheet) ---
$wsProjects = $wb.Worksheets.Item("projects_bin2")
$wsProjects.Range("A3").Value = "533valeo"
$wsProjects.Range("B3").Value = "192.168.277.27"
$wsProjects.Range("D3").ClearContents()

# --- 4. Settings_recources ---
$wsRecources = $wb.Worksheets.Item("Settings_recources")
$wsRecources.Range("A33").Value = "Nastavení spouštění TRIMAZKON v nabídce system tray"
$wsRecources.Range("B33").Value = "ano"
